# "delete 얀센 from 2nd dose"
#
# The 2nd-dose sheet had Janssen (얀센, a single-dose vaccine) counts mixed
# into several weekly cumulative totals. This script removes the Janssen
# contribution from each week's running total on the "2nd dose" sheet
# (rows 11-19, columns B-H), and fixes a fat-fingered value on the
# "1st dose" sheet (G12).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1st dose")
$ws2 = $wb.Worksheets.Item("2nd dose")

# --- 1st dose: fix mistyped value in G12 (33118792 -> 3318792) ---
$ws1.Range("G12").Value = 3318792

# --- 2nd dose: subtract the Janssen (얀센) amounts baked into the totals ---
# Each entry is written as a formula "<old total>-<janssen amount>" so the
# stored formula documents what was removed, matching how rows 14/16
# already recorded their Janssen add-ins as formulas.

$janssenSubtract = @{
    11 = @{ B = 17;     C = 802874; D = 197434; E = 60318;  F = 60206;  G = 8270;  H = 320 }
    12 = @{ B = 18;     C = 802951; D = 197528; E = 60379;  F = 60218;  G = 8272;  H = 320 }
    13 = @{ B = 18;     C = 802956; D = 197534; E = 60381;  F = 60222;  G = 8274;  H = 320 }
    15 = @{ B = 18;     C = 802983; D = 197544; E = 60393;  F = 60227;  G = 8276;  H = 321 }
    17 = @{ B = 18;     C = 802987; D = 197547; E = 60395;  F = 60231;  G = 8278;  H = 321 }
    18 = @{ B = 19;     C = 802987; D = 197549; E = 60396;  F = 60234;  G = 8278;  H = 321 }
    19 = @{ B = 22;     C = 824599; D = 218975; E = 65430;  F = 61377;  G = 8461;  H = 557 }
}

foreach ($row in $janssenSubtract.Keys) {
    foreach ($col in $janssenSubtract[$row].Keys) {
        $cell = $ws2.Range("$col$row")
        $base = [int64]$cell.Value2
        $sub = $janssenSubtract[$row][$col]
        $cell.Formula = "=$base-$sub"
    }
}

# Rows 14 and 16 already carried the Janssen amount as an explicit "+N" add-in
# formula (e.g. "=754939+18"); the edit simply drops that add-in, leaving the
# bare base formula behind.
$janssenRemoveAddin = @{
    14 = @{ B = 754939; C = 335462; D = 437733; E = 551208; F = 387426; G = 1517041; H = 1729588 }
    16 = @{ B = 883966; C = 454249; D = 641635; E = 693297; F = 545683; G = 1558023; H = 1750721 }
}

foreach ($row in $janssenRemoveAddin.Keys) {
    foreach ($col in $janssenRemoveAddin[$row].Keys) {
        $base = $janssenRemoveAddin[$row][$col]
        $ws2.Range("$col$row").Formula = "=$base"
    }
}

# --- Sheet view / active tab: selection moved to 2nd dose (B3 on 1st dose,
#     E11 on 2nd dose, which becomes the active sheet) ---
$ws1.Range("B3").Select()
$ws2.Activate()
$ws2.Range("E11").Select()
